# Round the emmean (B), SE (C), lower.CL (E) and upper.CL (F) columns
# to 3 decimal places for all data rows (2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "E", "F")

for ($row = 2; $row -le 25; $row++) {
    foreach ($col in $columns) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = [Math]::Round([double]$cell.Value2, 3)
    }
}
